$wb = $excel.ActiveWorkbook

# The "IncorrectUserID" sheet holds the text shown when a user enters an
# incorrect user ID / password. Update the expected "incorrect password"
# text (B2) to the new copy.
$ws = $wb.Worksheets.Item("IncorrectUserID")
$ws.Range("B2").Value = "Incorrect user ID or password."
$ws.Range("B:B").ColumnWidth = 27

# Make this sheet the active one (it becomes the selected tab).
$ws.Activate()
